$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "EN_RELATION" worksheet right after the existing "RELATION" sheet.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "EN_RELATION"

# Copy the header formatting (style) from RELATION!A1:C1 onto the new sheet's A1:C1,
# then overwrite the copied text with the new English header values.
$ws1.Range("A1:C1").Copy() | Out-Null
$newSheet.Range("A1:C1").PasteSpecial(-4122) | Out-Null

$newSheet.Range("A1").Value = "Object"
$newSheet.Range("B1").Value = "Relation"
$newSheet.Range("C1").Value = "Subject"

# Restore selections: RELATION sheet selects A1:C1, EN_RELATION sheet selects D6
# and becomes the active (tab-selected) sheet.
$ws1.Range("A1:C1").Select() | Out-Null
$newSheet.Range("D6").Select() | Out-Null
